$d = $word.ActiveDocument

# Locate the paragraph that still needs to be kept: the one ending the
# bibliography entry for "Thomson Pioneira (2008)." Everything from the
# paragraph right after it through the "(c) 2020 ... Creative Commons
# Attribution" paragraph (inclusive) must be removed, while the trailing
# empty paragraph and the page-break paragraph at the very end of the
# document must be preserved.

$count = $d.Paragraphs.Count

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*Thomson Pioneira*") {
        $startPara = $i + 1
    }
    if ($txt -like "*Powered by Jekyll*") {
        $endPara = $i
    }
}

if ($startPara -ne $null -and $endPara -ne $null -and $startPara -le $endPara) {
    $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}

$d.Saved = $false
